$wb = $excel.ActiveWorkbook
$src = $wb.Worksheets.Item("2025-08-28")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "2025-08-29"

# Copy header formatting (bold, border, center/top alignment) from the prior day sheet
$src.Range("A1:D1").Copy()
$newSheet.Range("A1:D1").PasteSpecial(-4122)
$newSheet.Application.CutCopyMode = $false

$newSheet.Cells.Item(1,1).Value = 'rank'
$newSheet.Cells.Item(1,2).Value = 'title'
$newSheet.Cells.Item(1,3).Value = 'author'
$newSheet.Cells.Item(1,4).Value = 'latest_episode'

$newSheet.Cells.Item(2,1).Value = 1
$newSheet.Cells.Item(2,2).Value = '「ククク……。奴は四天王の中でも最弱」と解雇された俺、なぜか勇者と聖女の師匠になる'
$newSheet.Cells.Item(2,3).Value = '漫画：芳橋アツシ 原作：延野正行 キャラクター原案：坂野杏梨'
$newSheet.Cells.Item(2,4).Value = '第41話　奴は帰りたい（前編）'
$newSheet.Cells.Item(3,1).Value = 2
$newSheet.Cells.Item(3,2).Value = '落ちこぼれだった兄が実は最強 ～史上最強の勇者は転生し、学園で無自覚に無双する～'
$newSheet.Cells.Item(3,3).Value = '村上よしゆき 茨木野 あるてら'
$newSheet.Cells.Item(3,4).Value = '第４１話　勇者、人魚王国を救い、歓迎される。あと、六邪神将が、全員来る（５）'
$newSheet.Cells.Item(4,1).Value = 3
$newSheet.Cells.Item(4,2).Value = 'バーサス'
$newSheet.Cells.Item(4,3).Value = '原作：ONE 漫画：あずま京太郎 構成：bose'
$newSheet.Cells.Item(4,4).Value = '第27話 幸せの在り処（2）'
$newSheet.Cells.Item(5,1).Value = 4
$newSheet.Cells.Item(5,2).Value = '最弱貴族に転生したので悪役たちを集めてみた'
$newSheet.Cells.Item(5,3).Value = '空野進 sorani ファルまろ'
$newSheet.Cells.Item(5,4).Value = '第11話　最弱貴族、背中を流してもらう（４）'
$newSheet.Cells.Item(6,1).Value = 5
$newSheet.Cells.Item(6,2).Value = 'デスゲームに巻き込まれた山本さん、気ままにゲームバランスを崩壊させる'
$newSheet.Cells.Item(6,3).Value = 'ぽち(原作) カモトタツヤ(作画) 久賀フーナ(キャラクター原案)'
$newSheet.Cells.Item(6,4).Value = '第6話（前編）'
$newSheet.Cells.Item(7,1).Value = 6
$newSheet.Cells.Item(7,2).Value = '転生したらスライムだった件 異聞 ～魔国暮らしのトリニティ～'
$newSheet.Cells.Item(7,3).Value = '伏瀬 戸野タエ みっつばー'
$newSheet.Cells.Item(7,4).Value = '第109話　開国祭開幕［その4］'
$newSheet.Cells.Item(8,1).Value = 7
$newSheet.Cells.Item(8,2).Value = '手に入れた催眠アプリで夢のハーレム生活を送りたい'
$newSheet.Cells.Item(8,3).Value = '櫻ヨル(漫画) みょん(原作) マッパニナッタ(キャラクター原案)'
$newSheet.Cells.Item(8,4).Value = '第11話①'
$newSheet.Cells.Item(9,1).Value = 8
$newSheet.Cells.Item(9,2).Value = '無職の英雄　別にスキルなんか要らなかったんだが'
$newSheet.Cells.Item(9,3).Value = '原作：九頭七尾・上田夢人 漫画：名苗秋緒'
$newSheet.Cells.Item(9,4).Value = '第51話'
$newSheet.Cells.Item(10,1).Value = 9
$newSheet.Cells.Item(10,2).Value = 'デスマーチからはじまる異世界狂想曲'
$newSheet.Cells.Item(10,3).Value = 'あやめぐむ(作画) 愛七ひろ(原作) ｓｈｒｉ(キャラクター原案)'
$newSheet.Cells.Item(10,4).Value = '第128話'
$newSheet.Cells.Item(11,1).Value = 10
$newSheet.Cells.Item(11,2).Value = 'うちの清楚系委員長がかつて中二病アイドルだったことを俺だけが知っている。'
$newSheet.Cells.Item(11,3).Value = '三上こた こばやし少女 寝子空兄 ゆがー'
$newSheet.Cells.Item(11,4).Value = '第3話　変わらないもの'
$newSheet.Cells.Item(12,1).Value = 11
$newSheet.Cells.Item(12,2).Value = '幼馴染のS級パーティーから追放された聖獣使い。万能支援魔法と仲間を増やして最強へ！'
$newSheet.Cells.Item(12,3).Value = '黒田高祥(作画) かなりつ(原作) 転(キャラクター原案)'
$newSheet.Cells.Item(12,4).Value = '第52話-2'
$newSheet.Cells.Item(13,1).Value = 12
$newSheet.Cells.Item(13,2).Value = '貧乏騎士に嫁入りしたはずが!? 〜野人令嬢は皇太子妃になっても竜を狩りたい〜'
$newSheet.Cells.Item(13,3).Value = '漫画：夏川そぞろ 原作：宮前葵 キャラクター原案：ののまろ'
$newSheet.Cells.Item(13,4).Value = '第12話①立太子式'
$newSheet.Cells.Item(14,1).Value = 13
$newSheet.Cells.Item(14,2).Value = '転生大聖女の異世界のんびり紀行'
$newSheet.Cells.Item(14,3).Value = 'キダニエル 四葉夕ト'
$newSheet.Cells.Item(14,4).Value = 'おまけ㊻	それって前世ブラックです！？'
$newSheet.Cells.Item(15,1).Value = 14
$newSheet.Cells.Item(15,2).Value = 'となりの席のヤツがそういう目で見てくる'
$newSheet.Cells.Item(15,3).Value = 'mmk'
$newSheet.Cells.Item(15,4).Value = '第47話 おやすみ'
$newSheet.Cells.Item(16,1).Value = 15
$newSheet.Cells.Item(16,2).Value = '魔術師クノンは見えている'
$newSheet.Cells.Item(16,3).Value = 'La-na(作画) 南野海風(原作) Ｌａｒｕｈａ(キャラクター原案)'
$newSheet.Cells.Item(16,4).Value = '第40話①'
$newSheet.Cells.Item(17,1).Value = 16
$newSheet.Cells.Item(17,2).Value = '【悲報】清楚系で売っていた底辺配信者、うっかり配信を切り忘れたままSS級モンスターを拳で殴り飛ばしてしまう'
$newSheet.Cells.Item(17,3).Value = 'アトハ NEO草野 pupps'
$newSheet.Cells.Item(17,4).Value = '第６話　【悲報】チェンジ作戦、大成功！？（３）'
$newSheet.Cells.Item(18,1).Value = 17
$newSheet.Cells.Item(18,2).Value = '生徒会にも穴はある！'
$newSheet.Cells.Item(18,3).Value = 'むちまろ'
$newSheet.Cells.Item(18,4).Value = '第134話	太賀のトラウマ'
$newSheet.Cells.Item(19,1).Value = 18
$newSheet.Cells.Item(19,2).Value = '異世界ゆるっとサバイバル生活～学校の皆と異世界の無人島に転移したけど俺だけ楽勝です～'
$newSheet.Cells.Item(19,3).Value = '西尾洋一(作画) 絢乃(原作) 乾和音(キャラクター原案) 株式会社一二三書房(監修)'
$newSheet.Cells.Item(19,4).Value = '第50話-2'
$newSheet.Cells.Item(20,1).Value = 19
$newSheet.Cells.Item(20,2).Value = '男女比1：5の世界でも普通に生きられると思った？　～激重感情な彼女たちが無自覚男子に翻弄されたら～'
$newSheet.Cells.Item(20,3).Value = '三藤 孝太郎(原作) 桃季憂(漫画) jimmy(キャラクター原案)'
$newSheet.Cells.Item(20,4).Value = '第10話-2'
$newSheet.Cells.Item(21,1).Value = 20
$newSheet.Cells.Item(21,2).Value = '町人Aは悪役令嬢をどうしても救いたい　～どぶと空と氷の姫君～'
$newSheet.Cells.Item(21,3).Value = '原作：一色孝太郎・Parum 漫画：目黒三吉'
$newSheet.Cells.Item(21,4).Value = '第39話 救出'
$newSheet.Cells.Item(22,1).Value = 21
$newSheet.Cells.Item(22,2).Value = 'まんきつしたい常連さん'
$newSheet.Cells.Item(22,3).Value = 'しんみりん(著者)'
$newSheet.Cells.Item(22,4).Value = '第47話前編'
$newSheet.Cells.Item(23,1).Value = 22
$newSheet.Cells.Item(23,2).Value = '放課後はケンカ最強のギャルに連れこまれる生活 彼女たちに好かれて、僕も最強に!?'
$newSheet.Cells.Item(23,3).Value = '亜逸(原作) あおやぎ孝夫(作画) ｋａｋａｏ(キャラクター原案)'
$newSheet.Cells.Item(23,4).Value = '第17話'
$newSheet.Cells.Item(24,1).Value = 23
$newSheet.Cells.Item(24,2).Value = '放課後、ファミレスで、クラスのあの子と。'
$newSheet.Cells.Item(24,3).Value = '麦子(漫画) 左リュウ(原作) magako(キャラクター原案)'
$newSheet.Cells.Item(24,4).Value = '第11話①'
$newSheet.Cells.Item(25,1).Value = 24
$newSheet.Cells.Item(25,2).Value = '規格外のダンジョン攻略者、実は異世界帰りの元勇者'
$newSheet.Cells.Item(25,3).Value = '作画：やまざき君 原作：榊与一'
$newSheet.Cells.Item(25,4).Value = '第6話(2)'
$newSheet.Cells.Item(26,1).Value = 25
$newSheet.Cells.Item(26,2).Value = '村人ですが何か？'
$newSheet.Cells.Item(26,3).Value = '鯖夢(作画) 白石新(原案・監修) 白蘇ふぁみ(キャラクター原案)'
$newSheet.Cells.Item(26,4).Value = '第92話'
$newSheet.Cells.Item(27,1).Value = 26
$newSheet.Cells.Item(27,2).Value = '継母の連れ子が元カノだった'
$newSheet.Cells.Item(27,3).Value = '草壁レイ(作画) 紙城境介(原作) たかやKi(キャラクター原案)'
$newSheet.Cells.Item(27,4).Value = '第62話前半'
$newSheet.Cells.Item(28,1).Value = 27
$newSheet.Cells.Item(28,2).Value = '実は俺、最強でした？'
$newSheet.Cells.Item(28,3).Value = '原作：澄守 彩 漫画：高橋 愛'
$newSheet.Cells.Item(28,4).Value = '第123話　王妃とハルト・後編'
$newSheet.Cells.Item(29,1).Value = 28
$newSheet.Cells.Item(29,2).Value = 'ダンジョンの幼なじみ'
$newSheet.Cells.Item(29,3).Value = '久真やすひさ(著者)'
$newSheet.Cells.Item(29,4).Value = '第56話'
$newSheet.Cells.Item(30,1).Value = 29
$newSheet.Cells.Item(30,2).Value = '「美人でお金持ちの彼女が欲しい」と言ったら、ワケあり女子がやってきた件。'
$newSheet.Cells.Item(30,3).Value = '白鷺六羽(作画) 小宮地千々(原作) Re岳(キャラクター原案) マイクロマガジン社(監修)'
$newSheet.Cells.Item(30,4).Value = '第28話-2'
$newSheet.Cells.Item(31,1).Value = 30
$newSheet.Cells.Item(31,2).Value = 'ギルド追放された雑用係の下剋上～超万能な生活スキルで世界最強～'
$newSheet.Cells.Item(31,3).Value = '原作／夜桜ユノ 漫画／柳輪 ネーム構成／ユーキあきら'
$newSheet.Cells.Item(31,4).Value = '第65話'
$newSheet.Cells.Item(32,1).Value = 31
$newSheet.Cells.Item(32,2).Value = '斎藤義龍に生まれ変わったので、織田信長に国譲りして長生きするのを目指します！'
$newSheet.Cells.Item(32,3).Value = '巽未頼 田村ゆうき マキムラシュンスケ'
$newSheet.Cells.Item(32,4).Value = '第74話「休暇」'
$newSheet.Cells.Item(33,1).Value = 32
$newSheet.Cells.Item(33,2).Value = 'ルパン三世 異世界の姫君（ネイバーワールドプリンセス）'
$newSheet.Cells.Item(33,3).Value = 'モンキー・パンチ／エム・ピー・ワークス 内々けやき 佐伯庸介 白狼'
$newSheet.Cells.Item(33,4).Value = '第107話：次元の優しいスナイプ'
$newSheet.Cells.Item(34,1).Value = 33
$newSheet.Cells.Item(34,2).Value = '異世界魔王と召喚少女の奴隷魔術'
$newSheet.Cells.Item(34,3).Value = '原作：むらさきゆきや 漫画：福田直叶 キャラクター原案：鶴崎貴大'
$newSheet.Cells.Item(34,4).Value = '『異世界魔王』ヒロイン総選挙 結果発表'
$newSheet.Cells.Item(35,1).Value = 34
$newSheet.Cells.Item(35,2).Value = '殺されたらゾンビになったので、進化しまくって無双しようと思います'
$newSheet.Cells.Item(35,3).Value = '漫画：朝ケ夜 原作：幸運ピエロ キャラクター原案：東西'
$newSheet.Cells.Item(35,4).Value = '第17話（前半） クラン「星覇」と序列戦①'
$newSheet.Cells.Item(36,1).Value = 35
$newSheet.Cells.Item(36,2).Value = '没落予定なので、鍛冶職人を目指す'
$newSheet.Cells.Item(36,3).Value = '石田彩(作画) CK(原作) かわく(キャラクター原案)'
$newSheet.Cells.Item(36,4).Value = '第101話'
$newSheet.Cells.Item(37,1).Value = 36
$newSheet.Cells.Item(37,2).Value = '転生したらスライムだった件　美食伝 ～ペコとリムルの料理手帖～'
$newSheet.Cells.Item(37,3).Value = '原作：伏瀬 漫画：中谷チカ キャラクター原案：みっつばー'
$newSheet.Cells.Item(37,4).Value = '第２４皿　魅惑の夜食‼ 極旨カツサンド‼（１）'
$newSheet.Cells.Item(38,1).Value = 37
$newSheet.Cells.Item(38,2).Value = '僕のカノジョ先生'
$newSheet.Cells.Item(38,3).Value = '星河蟹(作画) 孟倫（ＳＤｗｉｎｇ）(構成) 鏡遊(原作) おりょう(キャラクター原案)'
$newSheet.Cells.Item(38,4).Value = '76時間目-2'
$newSheet.Cells.Item(39,1).Value = 38
$newSheet.Cells.Item(39,2).Value = 'ダンジョンバンド'
$newSheet.Cells.Item(39,3).Value = '新挑限(著者)'
$newSheet.Cells.Item(39,4).Value = '♯８光の雨が降る夜に ②'
$newSheet.Cells.Item(40,1).Value = 39
$newSheet.Cells.Item(40,2).Value = '独身貴族は異世界を謳歌する ～結婚しない男の優雅なおひとりさまライフ～'
$newSheet.Cells.Item(40,3).Value = '漫画：駒鳥 ひわ 原作：錬金王 キャラクター原案：三登 いつき'
$newSheet.Cells.Item(40,4).Value = '第33話 独身貴族は見積もりを誤る（2）'
$newSheet.Cells.Item(41,1).Value = 40
$newSheet.Cells.Item(41,2).Value = 'ポンコツ勇者パーティー、竜をひろう'
$newSheet.Cells.Item(41,3).Value = '優風(著者)'
$newSheet.Cells.Item(41,4).Value = '第4話'
$newSheet.Cells.Item(42,1).Value = 41
$newSheet.Cells.Item(42,2).Value = 'モブだけど最強を目指します！　～ゲーム世界に転生した俺は自由に強さを追い求める～'
$newSheet.Cells.Item(42,3).Value = '反面教師(原作) 五條さやか(作画) 大熊猫介(キャラクター原案)'
$newSheet.Cells.Item(42,4).Value = '第13話'
$newSheet.Cells.Item(43,1).Value = 42
$newSheet.Cells.Item(43,2).Value = 'Only Sense Online ‐オンリーセンス・オンライン‐'
$newSheet.Cells.Item(43,3).Value = '羽仁倉雲(作画) アロハ座長(原作) ゆきさん(キャラクター原案)'
$newSheet.Cells.Item(43,4).Value = '第129話-2'
$newSheet.Cells.Item(44,1).Value = 43
$newSheet.Cells.Item(44,2).Value = '直径3cmの召喚陣<リミットリング>で「雑魚すら呼べない」と蔑まれた底辺召喚士が頂点に立つまで'
$newSheet.Cells.Item(44,3).Value = '作画：まっつー 原作：空松蓮司'
$newSheet.Cells.Item(44,4).Value = '第6話(2)'
$newSheet.Cells.Item(45,1).Value = 44
$newSheet.Cells.Item(45,2).Value = 'おねえさんと猫を飼う'
$newSheet.Cells.Item(45,3).Value = '上杉響士郎(著者)'
$newSheet.Cells.Item(45,4).Value = '第5話：おねえさんと猫の重さ'
$newSheet.Cells.Item(46,1).Value = 45
$newSheet.Cells.Item(46,2).Value = '君の刀が折れるまで ~月宮まつりの恋難き~'
$newSheet.Cells.Item(46,3).Value = 'イノウエ'
$newSheet.Cells.Item(46,4).Value = '第41話 決意'
$newSheet.Cells.Item(47,1).Value = 46
$newSheet.Cells.Item(47,2).Value = '北斗の拳 世紀末ドラマ撮影伝'
$newSheet.Cells.Item(47,3).Value = '原案/武論尊・原哲夫 漫画/倉尾宏'
$newSheet.Cells.Item(47,4).Value = '第75話 宙を舞う悪役俳優‼︎'
$newSheet.Cells.Item(48,1).Value = 47
$newSheet.Cells.Item(48,2).Value = '時間停止勇者―余命３日の設定じゃ世界を救うには短すぎる―'
$newSheet.Cells.Item(48,3).Value = '光永康則'
$newSheet.Cells.Item(48,4).Value = '第６８話『施錠停止』④'
$newSheet.Cells.Item(49,1).Value = 48
$newSheet.Cells.Item(49,2).Value = 'スキル【万物支配】に目覚めたおっさんは、ダンジョンで生計を立てることにしました～無職から始める支配者無双～'
$newSheet.Cells.Item(49,3).Value = '岸本和葉 原田 臙 シミズヒロノリ 吉武'
$newSheet.Cells.Item(49,4).Value = '第5話　一撃(後編)'
$newSheet.Cells.Item(50,1).Value = 49
$newSheet.Cells.Item(50,2).Value = '配信に致命的に向いていない女の子が迷宮で黙々と人助けする配信'
$newSheet.Cells.Item(50,3).Value = '下田将也(漫画) 佐藤悪糖(原作) 福きつね(キャラクター原案)'
$newSheet.Cells.Item(50,4).Value = '第2話前編'
$newSheet.Cells.Item(51,1).Value = 50
$newSheet.Cells.Item(51,2).Value = 'カナン様はあくまでチョロい'
$newSheet.Cells.Item(51,3).Value = 'nonco'
$newSheet.Cells.Item(51,4).Value = '第148話	カナンの布団の中'

# Restore original active sheet/selection so the workbook-level view state
# (active tab) matches the pre-edit workbook rather than the newly added sheet.
$wb.Worksheets.Item(1).Activate()
